$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.2377347397804473
$ws.Cells.Item(3, 2).Value = 0.1543604700559567
$ws.Cells.Item(4, 2).Value = 0.251585782613469
$ws.Cells.Item(5, 2).Value = 0.1494280070122514
$ws.Cells.Item(6, 2).Value = 0.1354613017779132
$ws.Cells.Item(7, 2).Value = 0.2650730800634213
$ws.Cells.Item(8, 2).Value = 0.2344869295150573
$ws.Cells.Item(9, 2).Value = 0.1464603184485922
$ws.Cells.Item(10, 2).Value = 0.2113211214126457
$ws.Cells.Item(11, 2).Value = 0.2075638418632239
$ws.Cells.Item(12, 2).Value = 0.1665166353100112
$ws.Cells.Item(13, 2).Value = 0.1579540858194672
$ws.Cells.Item(14, 2).Value = 0.1652945669008024
$ws.Cells.Item(15, 2).Value = 0.1645562403779922
$ws.Cells.Item(16, 2).Value = 0.182182048665788
$ws.Cells.Item(17, 2).Value = 0.2219619050606416
$ws.Cells.Item(18, 2).Value = 0.09117381186172682
$ws.Cells.Item(19, 2).Value = 0.2632210332655361
$ws.Cells.Item(20, 2).Value = 0.1820906027470262
$ws.Cells.Item(21, 2).Value = 0.1571869149263968
$ws.Cells.Item(22, 2).Value = 0.1666665446218413
$ws.Cells.Item(23, 2).Value = 0.1889785612249422
$ws.Cells.Item(24, 2).Value = 0.08767079386919532
$ws.Cells.Item(25, 2).Value = 0.1721018354523523
$ws.Cells.Item(26, 2).Value = 0.2742781765229103
$ws.Cells.Item(27, 2).Value = 0.1951429116059076
$ws.Cells.Item(28, 2).Value = 0.3070528500481298
$ws.Cells.Item(29, 2).Value = 0.1439333127970953
$ws.Cells.Item(30, 2).Value = 0.192540509710008
$ws.Cells.Item(31, 2).Value = 0.3022253590717894
$ws.Cells.Item(32, 2).Value = 0.1705504073475237
$ws.Cells.Item(33, 2).Value = 0.1842622732743352
$ws.Cells.Item(34, 2).Value = 0.2013297374803324
$ws.Cells.Item(35, 2).Value = 0.1924914849925276
$ws.Cells.Item(36, 2).Value = 0.1764049565438014
$ws.Cells.Item(37, 2).Value = 0.205713429138442
$ws.Cells.Item(38, 2).Value = 0.241736018927826
$ws.Cells.Item(39, 2).Value = 0.1511705503875471
$ws.Cells.Item(40, 2).Value = 0.2092311324332255
$ws.Cells.Item(41, 2).Value = 0.2027128115013695
$ws.Cells.Item(42, 2).Value = 0.1940841778297292
$ws.Cells.Item(43, 2).Value = 0.1546303521256976
$ws.Cells.Item(44, 2).Value = 0.1573236815058442
$ws.Cells.Item(45, 2).Value = 0.2093296792365324
$ws.Cells.Item(46, 2).Value = 0.1641695163468298
$ws.Cells.Item(47, 2).Value = 0.2452345544152963
$ws.Cells.Item(48, 2).Value = 0.214120442639659
$ws.Cells.Item(49, 2).Value = 0.2019740948408203
$ws.Cells.Item(50, 2).Value = 0.2275177740135915
$ws.Cells.Item(51, 2).Value = 0.1688688910828539
$ws.Cells.Item(52, 2).Value = 0.2583928921774317
$ws.Cells.Item(53, 2).Value = 0.1714816158709366
$ws.Cells.Item(54, 2).Value = 0.1876586990065328
$ws.Cells.Item(55, 2).Value = 0.2688365548518937
$ws.Cells.Item(56, 2).Value = 0.2286832407029697
$ws.Cells.Item(57, 2).Value = 0.118770108077467
$ws.Cells.Item(58, 2).Value = 0.2510735196905842
$ws.Cells.Item(59, 2).Value = 0.2313059547443149
$ws.Cells.Item(60, 2).Value = 0.224577544702168
$ws.Cells.Item(61, 2).Value = 0.2323060716593902
$ws.Cells.Item(62, 2).Value = 0.16973509223877
$ws.Cells.Item(63, 2).Value = 0.152090954373878
$ws.Cells.Item(64, 2).Value = 0.2400130043555971
$ws.Cells.Item(65, 2).Value = 0.1920021496176086
$ws.Cells.Item(66, 2).Value = 0.1881974001442445
$ws.Cells.Item(67, 2).Value = 0.218009627802776
$ws.Cells.Item(68, 2).Value = 0.1713604323410646
$ws.Cells.Item(69, 2).Value = 0.1934907850019723
$ws.Cells.Item(70, 2).Value = 0.1688317038968409
$ws.Cells.Item(71, 2).Value = 0.1697854750245221
$ws.Cells.Item(72, 2).Value = 0.2095644747452401
$ws.Cells.Item(73, 2).Value = 0.1944915149864047
$ws.Cells.Item(74, 2).Value = 0.16980769471393
$ws.Cells.Item(75, 2).Value = 0.2023392561536684
$ws.Cells.Item(76, 2).Value = 0.1830764900074947
$ws.Cells.Item(77, 2).Value = 0.1653809993627637
$ws.Cells.Item(78, 2).Value = 0.1633942294992086
